$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 11 new columns before column D (old D:K -> new O:V), preserving
#    the existing CardHolderName/CardCvvNo/CardNo/Address/City/Zipcode/
#    Phone_Number/BankName block (with its formatting) intact but shifted.
# ---------------------------------------------------------------------------
$ws.Range("D1:N1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 2. New header row (row 1) for the inserted guest-details columns.
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "Firstname_GuestOne"
$ws.Range("E1").Value = "Middlename_GuestOne"
$ws.Range("F1").Value = "Lastname_GuestOne"
$ws.Range("G1").Value = "Email"
$ws.Range("H1").Value = "Street_Address"
$ws.Range("I1").Value = "Cty"
$ws.Range("J1").Value = "Zpcode"
$ws.Range("K1").Value = "Phone"
$ws.Range("L1").Value = "Firstname_GuestTwo"
$ws.Range("M1").Value = "Middlename_GuestTwo"
$ws.Range("N1").Value = "Lastname_GuestTwo"

# ---------------------------------------------------------------------------
# 3. New data row (row 2) for the inserted guest-details columns.
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = "Monish"
$ws.Range("E2").Value = "N"
$ws.Range("F2").Value = "Luthra"
$ws.Range("G2").Value = "ajit_nakum@odysseussolutions.com"
$ws.Range("H2").Value = "White house"
$ws.Range("I2").Value = "Miami"
$ws.Range("J2").Value = 10245
$ws.Range("K2").Value = 4565289563
$ws.Range("L2").Value = "Anna"
$ws.Range("M2").Value = "M"
$ws.Range("N2").Value = "Luthra"

# Give the email cell an actual mailto hyperlink + the Hyperlink style
# (re-applied after Add so the cell keeps reusing the workbook's existing
# "Hyperlink" cell style instead of acquiring a throwaway duplicate one).
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:ajit_nakum@odysseussolutions.com")
$ws.Range("G2").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 4. Update the (shifted) card-number value that now lives in O2.
# ---------------------------------------------------------------------------
$ws.Range("O2").NumberFormat = "@"
$ws.Range("O2").Value = "4387751111111111"

# ---------------------------------------------------------------------------
# 5. Column widths for the whole, now-wider sheet (A:V).
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 19.166666666666668    # D  width 20
$ws.Columns.Item(5).ColumnWidth = 21.666666666666668    # E  width 22.5703125
$ws.Columns.Item(6).ColumnWidth = 18.666666666666668    # F  width 19.5703125
$ws.Columns.Item(7).ColumnWidth = 33.333333333333336    # G  width 34.140625
$ws.Columns.Item(8).ColumnWidth = 13.833333333333334    # H  width 14.7109375
$ws.Range("I1:K1").EntireColumn.ColumnWidth = 13.833333333333334  # I:K width 14.7109375
$ws.Columns.Item(12).ColumnWidth = 19.166666666666668   # L  width 20
$ws.Columns.Item(13).ColumnWidth = 21.666666666666668   # M  width 22.5703125
$ws.Columns.Item(14).ColumnWidth = 19.166666666666668   # N  width 20

# ---------------------------------------------------------------------------
# 6. The sheet's used range/dimension extends two columns further right
#    (through X) than the populated data (which ends at V) -- touch X2 (with
#    a formatting no-op, so it stays otherwise blank/default-styled) so the
#    saved dimension/row-span bookkeeping matches.
# ---------------------------------------------------------------------------
$ws.Range("X2").Font.Bold = $false

# ---------------------------------------------------------------------------
# 7. Final selection, matching the saved view state.
# ---------------------------------------------------------------------------
[void]$ws.Range("M9").Select()

Write-Output "done"
